# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2 through 28
# from 2023-09-06 (serial 45175) to 2023-09-14 (serial 45183).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
